# Apply "dSF" (column F) corrections to stanek_ryne.xlsx
# Commit message: "repull data, push all data, mean calculation"
# These rows had their F (dSF) value re-pulled from source data and now
# differ from the original E (dS0) value they used to mirror.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 0
    20 = 0
    22 = 1
    26 = 0
    28 = -3
    29 = 3
    31 = 0
    32 = -2
    39 = 2
    42 = 1
    59 = 1
    60 = -1
    64 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
